# "mongo pass as env" - refresh the subject names on the Fourth Year
# timetable sheets and move the active view to the last sheet (mr).

$wb = $excel.ActiveWorkbook

# --- cs sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("cs")
$ws.Range("B2").Value = "Dc"
$ws.Range("B3").Value = "Embedded systems"
$ws.Range("B3").Select() | Out-Null

# --- ec sheet -----------------------------------------------------------
$ws = $wb.Worksheets.Item("ec")
$ws.Range("B2").Value = "Control Systems"
$ws.Range("B3").Value = "VLSI Design"
$ws.Range("B2:B3").Select() | Out-Null
$ws.PageSetup.Orientation = 1

# --- ce sheet -------------------------------------------------------------
$ws = $wb.Worksheets.Item("ce")
$ws.Range("B2").Value = "Hydraulics"
$ws.Range("B3").Value = "Geomatics"
$ws.Range("B2:B3").Select() | Out-Null

# --- ee sheet -------------------------------------------------------------
$ws = $wb.Worksheets.Item("ee")
$ws.Range("B2").Value = "Electromagnetic Theory"
$ws.Range("B3").Value = "Renewable Energy Sources"
$ws.Range("C9").Select() | Out-Null

# --- me sheet -------------------------------------------------------------
$ws = $wb.Worksheets.Item("me")
$ws.Range("B2").Value = "Theory of Machines"
$ws.Range("B3").Value = "Engineering Thermodynamics"
$ws.Range("B2:B3").Select() | Out-Null

# --- mr sheet (ends up the active tab / selected sheet) --------------------
$ws = $wb.Worksheets.Item("mr")
$ws.Range("B2").Value = "System Dynamics and Analysis"
$ws.Range("B3").Value = "Instrumentation and Measurement"
$ws.Columns.Item(2).ColumnWidth = 31.8
$ws.Activate()
$ws.Range("D5").Select() | Out-Null
